$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 83, pushing the
# previously existing rows 83-200 down to rows 84-201.
$ws.Rows(83).Insert()

$ws.Range("A83").Value = 7
$ws.Range("B83").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C83").Value = "Ñuble"
$ws.Range("D83").Value2 = 44579
$ws.Range("E83").Value = 16
$ws.Range("F83").Value = 100112009
$ws.Range("G83").Value = "Acelga"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 120
$ws.Range("K83").Value = 350
$ws.Range("L83").Value = 400
$ws.Range("M83").Value = 375
$ws.Range("N83").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O83").Value = "Provincia de Diguillín"
$ws.Range("P83").Value = 375
$ws.Range("Q83").Value = 1
$ws.Range("R83").Value = "Hortaliza"
